$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Clinical only, 0-3 mo.
$ws.Range("C2").Value = 0.286
$ws.Range("D2").Value = 20.721
$ws.Range("E2").Value = 14.936
$ws.Range("F2").Value = 231.748
$ws.Range("G2").Value = 6

# Row 3: Clinical only, 3-6 mo.
$ws.Range("C3").Value = -1.007
$ws.Range("D3").Value = 30.081
$ws.Range("E3").Value = 25.217
$ws.Range("F3").Value = 233.7
$ws.Range("G3").Value = 6

# Row 4: Clinical only, 6-12 mo.
$ws.Range("C4").Value = -1.212
$ws.Range("D4").Value = 37.388
$ws.Range("E4").Value = 29.495
$ws.Range("F4").Value = 238.376
$ws.Range("G4").Value = 6

# Row 5: Clinical + wb-FA, 0-3 mo.
$ws.Range("C5").Value = 0.074
$ws.Range("D5").Value = 23.599
$ws.Range("E5").Value = 17.587
$ws.Range("F5").Value = 233.555
$ws.Range("G5").Value = 6

# Row 6: Clinical + wb-FA, 3-6 mo.
$ws.Range("C6").Value = -1.195
$ws.Range("D6").Value = 31.459
$ws.Range("E6").Value = 25.254
$ws.Range("F6").Value = 233.852
$ws.Range("G6").Value = 6

# Row 7: Clinical + wb-FA, 6-12 mo.
$ws.Range("C7").Value = -1.351
$ws.Range("D7").Value = 38.548
$ws.Range("E7").Value = 29.6
$ws.Range("F7").Value = 240.313
$ws.Range("G7").Value = 6

# Rows 8-10: Clinical + tractFA - G column becomes an empty text cell
# (was numeric 0) instead of being cleared entirely. A leading apostrophe
# forces Excel to store it as literal (empty) text; resetting the style
# afterwards drops the transient quote-prefix formatting it introduces.
$ws.Range("G8").Value = "'"
$ws.Range("G8").Style = "Normal"
$ws.Range("G9").Value = "'"
$ws.Range("G9").Style = "Normal"
$ws.Range("G10").Value = "'"
$ws.Range("G10").Style = "Normal"
